$d = $word.ActiveDocument

# Find the end of "... diagnostiziert." (last run of that paragraph)
# and append the new sentence about the therapist's app suggestion as
# two new runs (same Candara font), without touching existing runs.
$r = $d.Content
$r.Find.Execute("diagnostiziert.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

$insertStart = $r.Start

$part1 = " Seine Therapeutin schlug ihm vor die Applikation „E-Mood-Tracker“ zu verwenden, damit er seine tägliche Stimmung leichter aufzeichnen kann, ohne stän"
$part2 = "dig zur Klinik fahren zu müssen und sie ihm von der Klinik aus seine Medikamente regeln kann."

$r.InsertAfter($part1)
$len1 = $part1.Length
$run1 = $d.Range($insertStart, $insertStart + $len1)
$run1.Font.Name = "Candara"

$r2 = $d.Range($insertStart + $len1, $insertStart + $len1)
$r2.InsertAfter($part2)
$len2 = $part2.Length
$run2 = $d.Range($insertStart + $len1, $insertStart + $len1 + $len2)
$run2.Font.Name = "Candara"
